# Fill in the previously-blank C11:I11 cells (GDP, UEMP, CPI, LTRate,
# EURUSD, WTI, RPP for the 2017 Q4 row) with computed scenario values on
# both "Test 1" and "Test 2" worksheets.

$wb = $excel.ActiveWorkbook

$values = @(0.7935213953370379, -0.39999999999999947, 0.2989451731422861, -0.051999999999999935, 1.5829618029997903, 16.12947350163202, 1.6798418972332)

foreach ($sheetName in @("Test 1", "Test 2")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($i = 0; $i -lt $values.Length; $i++) {
        # Columns C..I correspond to index 3..9
        $col = 3 + $i
        $ws.Cells.Item(11, $col).Value = $values[$i]
    }
}
